$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-converted to numbers by Excel (they must remain text, matching
# the original inlineStr cell type).
$textCells = @("D5","D6","D9","D11","D13","D19","D20","D21","D22","D24","D26","D27","D28","D31","D32","D33","D35","D40","D41","D43","D44","D45","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range("D2").Value = "67.378.83"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "3.772.63"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "591.63"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("D6").Value = "171.32"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("D7").Value = "3.770.77"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").Value = "6.23"
$ws.Range("E11").Value = "  -5.23%  "
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("D13").Value = "37.67"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").Value = "4.397.56"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "3.772.38"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "67.466.35"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  -4.88%  "
$ws.Range("D20").Value = "15.98"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "485.66"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "9.14"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "83.87"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("E25").Value = "  -9.97%  "
$ws.Range("D26").Value = "0.0000140"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").Value = "12.13"
$ws.Range("E27").Value = "  -5.85%  "
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  -11.79%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "2.38"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").Value = "32.32"
$ws.Range("E32").Value = "  +6.75%  "
$ws.Range("D33").Value = "7.75"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("E39").Value = "  -8.26%  "
$ws.Range("D40").Value = "448.82"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").Value = "48.74"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E42").Value = "  -3.60%  "
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  -7.25%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "41.40"
$ws.Range("E44").Value = "  -9.57%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "8.23"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").Value = "2.824.79"
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").Value = "140.19"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("D50").Value = "25.75"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("D51").Value = "23.11"
$ws.Range("E51").Value = "  +7.41%  "
